$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.090.52'
$ws.Range("E2").Value = '  -2.11%  '

$ws.Range("D3").Value = '''3.423.30'
$ws.Range("E3").Value = '  -2.39%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''579.92'
$ws.Range("E5").Value = '  -1.76%  '

$ws.Range("D6").Value = '''128.93'
$ws.Range("E6").Value = '  -4.18%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '''0.479'
$ws.Range("E8").Value = '  -1.80%  '

$ws.Range("D9").Value = '''7.55'
$ws.Range("E9").Value = '  +1.60%  '

$ws.Range("D10").Value = '''0.123'
$ws.Range("E10").Value = '  -0.98%  '

$ws.Range("D11").Value = '''0.381'
$ws.Range("E11").Value = '  -1.02%  '

$ws.Range("D12").Value = '''4.009.61'
$ws.Range("E12").Value = '  -2.43%  '

$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("D14").Value = '''0.0000176'
$ws.Range("E14").Value = '  -2.93%  '

$ws.Range("D15").Value = '''3.426.10'
$ws.Range("E15").Value = '  -2.34%  '

$ws.Range("D16").Value = '''63.120.73'
$ws.Range("E16").Value = '  -2.11%  '

$ws.Range("D17").Value = '''25.03'
$ws.Range("E17").Value = '  -2.70%  '

$ws.Range("D18").Value = '''9.79'
$ws.Range("E18").Value = '  -1.34%  '

$ws.Range("D19").Value = '''5.68'
$ws.Range("E19").Value = '  -1.69%  '

$ws.Range("D20").Value = '''13.25'
$ws.Range("E20").Value = '  -2.59%  '

$ws.Range("D21").Value = '''382.09'
$ws.Range("E21").Value = '  -3.03%  '

$ws.Range("D22").Value = '''0.565'
$ws.Range("E22").Value = '  -1.46%  '

$ws.Range("D23").Value = '''3.562.03'
$ws.Range("E23").Value = '  -2.38%  '

$ws.Range("D24").Value = '''73.29'
$ws.Range("E24").Value = '  -1.84%  '

$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").Value = '''0.0000110'
$ws.Range("E26").Value = '  -5.59%  '

$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("D28").Value = '''7.07'
$ws.Range("E28").Value = '  -4.56%  '

$ws.Range("D29").Value = '''2.19'
$ws.Range("E29").Value = '  -3.52%  '

$ws.Range("D30").Value = '''7.91'
$ws.Range("E30").Value = '  -4.23%  '

$ws.Range("D31").Value = '''0.153'
$ws.Range("E31").Value = '  -1.10%  '

$ws.Range("D32").Value = '''1.41'
$ws.Range("E32").Value = '  -5.73%  '

$ws.Range("D33").Value = '''3.452.98'
$ws.Range("E33").Value = '  -2.18%  '

$ws.Range("D35").Value = '''22.71'
$ws.Range("E35").Value = '  -3.37%  '

$ws.Range("D36").Value = '''5.30'
$ws.Range("E36").Value = '  -1.43%  '

$ws.Range("D37").Value = '''6.78'
$ws.Range("E37").Value = '  -2.13%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '''1.52'
$ws.Range("E38").Value = '  -2.42%  '

$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '''163.55'
$ws.Range("E39").Value = '  -2.48%  '

$ws.Range("D40").Value = '''0.0769'
$ws.Range("E40").Value = '  -2.32%  '

$ws.Range("D41").Value = '''0.782'
$ws.Range("E41").Value = '  -3.68%  '

$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("D43").Value = '''41.15'
$ws.Range("E43").Value = '  -1.43%  '

$ws.Range("D44").Value = '''4.32'
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").Value = '''1.60'
$ws.Range("E45").Value = '  -4.06%  '

$ws.Range("D46").Value = '''23.25'
$ws.Range("E46").Value = '  -7.47%  '

$ws.Range("E47").Value = '  -6.40%  '

$ws.Range("D48").Value = '''6.73'
$ws.Range("E48").Value = '  -0.84%  '

$ws.Range("D49").Value = '''0.884'
$ws.Range("E49").Value = '  -1.31%  '

$ws.Range("D50").Value = '''2.247.88'
$ws.Range("E50").Value = '  -6.89%  '

$ws.Range("D51").Value = '''0.0251'
$ws.Range("E51").Value = '  -3.19%  '
